$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "34.648.29"
$ws.Range("E2").Value = "  +1.41%  "
Set-TextValue "D3" "1.806.46"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.30%  "
Set-TextValue "D5" "227.76"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("E7").Value = "  -0.29%  "
Set-TextValue "D8" "32.93"
$ws.Range("E8").Value = "  +3.72%  "
$ws.Range("E9").Value = "  +1.84%  "
Set-TextValue "D10" "0.0698"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  +0.57%  "
Set-TextValue "D12" "2.063.90"
$ws.Range("E12").Value = "  +1.26%  "
Set-TextValue "D13" "11.21"
$ws.Range("E13").Value = "  +2.57%  "
Set-TextValue "D14" "1.797.15"
$ws.Range("E14").Value = "  +0.53%  "
Set-TextValue "D15" "0.641"
$ws.Range("E15").Value = "  +2.78%  "
Set-TextValue "D16" "34.622.06"
$ws.Range("E16").Value = "  +1.37%  "
Set-TextValue "D17" "4.35"
$ws.Range("E17").Value = "  +3.79%  "
Set-TextValue "D18" "69.05"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("E19").Value = "  +0.39%  "
Set-TextValue "D20" "247.98"
$ws.Range("E20").Value = "  +0.80%  "
Set-TextValue "D21" "11.37"
$ws.Range("E21").Value = "  +3.56%  "
$ws.Range("E22").Value = "  -0.25%  "
Set-TextValue "D23" "4.21"
$ws.Range("E23").Value = "  +2.83%  "
Set-TextValue "D24" "170.70"
$ws.Range("E24").Value = "  +5.08%  "
$ws.Range("E25").Value = "  +1.68%  "
Set-TextValue "D26" "7.37"
$ws.Range("E26").Value = "  +2.58%  "
Set-TextValue "D27" "16.70"
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("E30").Value = "  +11.48%  "
Set-TextValue "D31" "0.0528"
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "3.83"
$ws.Range("E32").Value = "  +2.47%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D33" "1.24"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("E34").Value = "  +2.98%  "
Set-TextValue "D35" "1.433.00"
$ws.Range("E35").Value = "  -0.57%  "
Set-TextValue "D36" "2.62"
$ws.Range("E36").Value = "  +8.96%  "
Set-TextValue "D37" "0.678"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("E39").Value = "  +0.82%  "
Set-TextValue "D40" "85.45"
$ws.Range("E40").Value = "  +6.51%  "
Set-TextValue "D41" "0.953"
$ws.Range("E41").Value = "  +3.08%  "
Set-TextValue "D42" "2.40"
$ws.Range("E42").Value = "  +1.80%  "
$ws.Range("E43").Value = "  +3.46%  "
Set-TextValue "D44" "13.93"
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("E45").Value = "  +2.83%  "
Set-TextValue "D46" "6.13"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("E47").Value = "  +0.65%  "
Set-TextValue "D48" "1.962.14"
$ws.Range("E48").Value = "  +1.13%  "
Set-TextValue "D49" "106.09"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  -4.41%  "
